$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 24394.16842127913
$ws.Range("D3").Value = 1877.883891196722
$ws.Range("E3").Value = 4112.394024069348

# Row 4 (std)
$ws.Range("B4").Value = 10563.3427216366
$ws.Range("D4").Value = 1154.239585269906
$ws.Range("E4").Value = 3102.312427237336

# Row 5 (min)
$ws.Range("B5").Value = 4319.034000000001
$ws.Range("D5").Value = 0.004
$ws.Range("E5").Value = 46.03400000000001

# Row 6 (25%)
$ws.Range("B6").Value = 15985.026
$ws.Range("D6").Value = 790.005
$ws.Range("E6").Value = 2021.002

# Row 7 (50%)
$ws.Range("B7").Value = 22126.56050000007
$ws.Range("D7").Value = 1845.01
$ws.Range("E7").Value = 3260.008

# Row 8 (75%)
$ws.Range("B8").Value = 32806.71175000008
$ws.Range("D8").Value = 3002.003
$ws.Range("E8").Value = 4996.013249999999

# Row 9 (max)
$ws.Range("B9").Value = 62322.24400000002
$ws.Range("D9").Value = 6053.007000000001
$ws.Range("E9").Value = 32035.006

# Row 10 (Total)
$ws.Range("F10").Value = 12821574922.22701

# Row 11 (Residential)
$ws.Range("G11").Value = 0.7544381176757862

# Row 12 (Community)
$ws.Range("F12").Value = 987015773.2130002
$ws.Range("G12").Value = 0.0769808529139385

# Row 13 (IGA)
$ws.Range("F13").Value = 2161474299.05
$ws.Range("G13").Value = 0.1685810294102753
